$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

Set-TextValue "D2" "243.03"
Set-TextValue "E2" "-0.72%"

Set-TextValue "D3" "29.97"
Set-TextValue "E3" "13.04%"

Set-TextValue "D4" "5.146"
Set-TextValue "E4" "0.40%"

Set-TextValue "D5" "0.05668"
Set-TextValue "E5" "1.44%"

Set-TextValue "D6" "6.520"
Set-TextValue "E6" "0.81%"

Set-TextValue "D7" "0.8392"
Set-TextValue "E7" "2.63%"

Set-TextValue "D8" "0.8596"
Set-TextValue "E8" "2.87%"

Set-TextValue "B9" "WazirX"
Set-TextValue "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1337"
Set-TextValue "E9" "0.16%"

Set-TextValue "B10" "MandalaExchangeToken"
Set-TextValue "C10" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.06910"
Set-TextValue "E10" "-0.97%"

Set-TextValue "B11" "BitrueCoin"
Set-TextValue "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.02861"
Set-TextValue "E11" "-0.86%"

Set-TextValue "B12" "BitMartToken"
Set-TextValue "C12" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D12" "0.09370"
Set-TextValue "E12" "-0.20%"

Set-TextValue "B13" "BitForexToken"
Set-TextValue "C13" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D13" "0.001520"
Set-TextValue "E13" "0.89%"

Set-TextValue "B14" "CoinExToken"
Set-TextValue "C14" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D14" "0.04156"
Set-TextValue "E14" "-10.33%"

Set-TextValue "B15" "One"
Set-TextValue "C15" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D15" "0.0006017"
Set-TextValue "E15" "-94.01%"

Set-TextValue "D16" "0.005998"
Set-TextValue "E16" "-3.21%"

Set-TextValue "D17" "3.508"
Set-TextValue "E17" "-3.84%"

Set-TextValue "D18" "3.018"
Set-TextValue "E18" "-0.50%"

Set-TextValue "D19" "2.127"
Set-TextValue "E19" "-2.57%"

Set-TextValue "E20" "1.26%"

Set-TextValue "D21" "0.03253"
Set-TextValue "E21" "4.10%"

Set-TextValue "E22" "-0.34%"

Set-TextValue "D23" "3.613"
Set-TextValue "E23" "-3.29%"

Set-TextValue "E24" "-0.10%"

Set-TextValue "D25" "0.001209"
Set-TextValue "E25" "-3.02%"

Set-TextValue "D26" "0.004451"
Set-TextValue "E26" "-1.10%"

Set-TextValue "E27" "22.80%"

Set-TextValue "D28" "0.0001397"
Set-TextValue "E28" "0.25%"

Set-TextValue "D40" "0.03711"
Set-TextValue "E40" "1.94%"

Set-TextValue "D41" "0.005326"
Set-TextValue "E41" "-13.84%"

Set-TextValue "D42" "0.1057"
Set-TextValue "E42" "0.65%"

Set-TextValue "D43" "0.002089"
Set-TextValue "E43" "-12.98%"

Set-TextValue "D44" "0.009738"
Set-TextValue "E44" "9.80%"

Set-TextValue "D45" "0.00005105"
Set-TextValue "E45" "-4.77%"

Set-TextValue "E46" "-0.10%"

Set-TextValue "D47" "0.09994"
Set-TextValue "E47" "-30.61%"

Set-TextValue "D48" "0.002700"
Set-TextValue "E48" "15.64%"

Set-TextValue "E49" "-0.10%"

Set-TextValue "E50" "-0.10%"
